# Adds 2022-Q3 data:
#  - duplicates the "2022-Q2" sheet (preserving its data as the new "2022-Q2" tab)
#  - turns the original "2022-Q2" sheet into the new "2022-Q3" tab with fresh figures
#  - inserts the corresponding summary row on "总计"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the "2022-Q2" worksheet so its current data is preserved under
#    the same tab name once the original sheet is repurposed for 2022-Q3.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)
$q2dup = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 2) Turn the original sheet (still named "2022-Q2", now at position 2) into
#    the "2022-Q3" sheet with the latest fund figures; only then can the
#    duplicate reclaim the "2022-Q2" name.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"
$q2dup.Name = "2022-Q2"

$q3.Range("B2").Value = "159941"
$q3.Range("C2").Value = "广发纳斯达克100ETF（QDII）"
$q3.Range("D2").Value = "106.15"
$q3.Range("E2").Value = "91.14"
$q3.Range("F2").Value = "3.22"
$q3.Range("G2").Value = "3.4180"
$q3.Range("H2").Value = 5

$q3.Range("B3").Value = "513100"
$q3.Range("C3").Value = "国泰纳斯达克100（QDII-ETF）"
$q3.Range("D3").Value = "46.54"
$q3.Range("E3").Value = "91.35"
$q3.Range("F3").Value = "3.31"
$q3.Range("G3").Value = "1.5405"
$q3.Range("H3").Value = 5

$q3.Range("B4").Value = "513500"
$q3.Range("C4").Value = "博时标普500ETF（QDII）"
$q3.Range("D4").Value = "71.37"
$q3.Range("E4").Value = "96.44"
$q3.Range("F4").Value = "1.65"
$q3.Range("G4").Value = "1.1776"
$q3.Range("H4").Value = 6

$q3.Range("B5").Value = "040047"
$q3.Range("C5").Value = "华安纳斯达克100指数（QDII）美元现钞A"
$q3.Range("D5").Value = "24.52"
$q3.Range("E5").Value = "92.09"
$q3.Range("F5").Value = "3.27"
$q3.Range("G5").Value = "0.8018"
$q3.Range("H5").Value = 5

$q3.Range("B6").Value = "040048"
$q3.Range("C6").Value = "华安纳斯达克100指数（QDII）美元现汇A"
$q3.Range("D6").Value = "24.52"
$q3.Range("E6").Value = "92.09"
$q3.Range("F6").Value = "3.27"
$q3.Range("G6").Value = "0.8018"
$q3.Range("H6").Value = 5

$q3.Range("B7").Value = "040046"
$q3.Range("C7").Value = "华安纳斯达克100指数（QDII）人民币A"
$q3.Range("D7").Value = "22.21"
$q3.Range("E7").Value = "92.09"
$q3.Range("F7").Value = "3.27"
$q3.Range("G7").Value = "0.7263"
$q3.Range("H7").Value = 5

$q3.Range("B8").Value = "160213"
$q3.Range("C8").Value = "国泰纳斯达克100指数（QDII）"
$q3.Range("D8").Value = "15.14"
$q3.Range("E8").Value = "85.81"
$q3.Range("F8").Value = "3.16"
$q3.Range("G8").Value = "0.4784"
$q3.Range("H8").Value = 5

$q3.Range("B9").Value = "000834"
$q3.Range("C9").Value = "大成纳斯达克100指数（QDII）"
$q3.Range("D9").Value = "14.15"
$q3.Range("E9").Value = "85.22"
$q3.Range("F9").Value = "3.06"
$q3.Range("G9").Value = "0.4330"
$q3.Range("H9").Value = 5

$q3.Range("B10").Value = "003722"
$q3.Range("C10").Value = "易方达纳斯达克100指数美元（QDII-LOF）A"
$q3.Range("D10").Value = "7.72"
$q3.Range("E10").Value = "90.67"
$q3.Range("F10").Value = "3.25"
$q3.Range("G10").Value = "0.2509"
$q3.Range("H10").Value = 5

$q3.Range("B11").Value = "161130"
$q3.Range("C11").Value = "易方达纳斯达克100指数人民币（QDII-LOF）"
$q3.Range("D11").Value = "7.72"
$q3.Range("E11").Value = "90.67"
$q3.Range("F11").Value = "3.25"
$q3.Range("G11").Value = "0.2509"
$q3.Range("H11").Value = 5

$q3.Range("B12").Value = "014978"
$q3.Range("C12").Value = "华安纳斯达克100指数（QDII）人民币C"
$q3.Range("D12").Value = "2.31"
$q3.Range("E12").Value = "92.09"
$q3.Range("F12").Value = "3.27"
$q3.Range("G12").Value = "0.0755"
$q3.Range("H12").Value = 5

$q3.Range("B13").Value = "161125"
$q3.Range("C13").Value = "易方达标普500指数（QDII-LOF）人民币"
$q3.Range("D13").Value = "4.74"
$q3.Range("E13").Value = "90.72"
$q3.Range("F13").Value = "1.56"
$q3.Range("G13").Value = "0.0739"
$q3.Range("H13").Value = 6

$q3.Range("B14").Value = "012860"
$q3.Range("C14").Value = "易方达标普500指数（QDII-LOF）人民币 C"
$q3.Range("D14").Value = "4.74"
$q3.Range("E14").Value = "90.72"
$q3.Range("F14").Value = "1.56"
$q3.Range("G14").Value = "0.0739"
$q3.Range("H14").Value = 6

$q3.Range("B15").Value = "003718"
$q3.Range("C15").Value = "易方达标普500指数（QDII-LOF）美元A"
$q3.Range("D15").Value = "4.66"
$q3.Range("E15").Value = "90.72"
$q3.Range("F15").Value = "1.56"
$q3.Range("G15").Value = "0.0727"
$q3.Range("H15").Value = 6

$q3.Range("B16").Value = "159632"
$q3.Range("C16").Value = "华安纳斯达克100ETF（QDII）"
$q3.Range("D16").Value = "1.51"
$q3.Range("E16").Value = "89.05"
$q3.Range("F16").Value = "3.17"
$q3.Range("G16").Value = "0.0479"
$q3.Range("H16").Value = 5

$q3.Range("B17").Value = "005698"
$q3.Range("C17").Value = "华夏全球科技先锋混合（QDII）"
$q3.Range("D17").Value = "0.59"
$q3.Range("E17").Value = "86.79"
$q3.Range("F17").Value = "6.46"
$q3.Range("G17").Value = "0.0381"
$q3.Range("H17").Value = 6

$q3.Range("B18").Value = "159612"
$q3.Range("C18").Value = "国泰标普500ETF（QDII）"
$q3.Range("D18").Value = "0.55"
$q3.Range("E18").Value = "91.40"
$q3.Range("F18").Value = "1.58"
$q3.Range("G18").Value = "0.0087"
$q3.Range("H18").Value = 6

$q3.Range("B19").Value = "006555"
$q3.Range("C19").Value = "浦银安盛全球智能科技股票（QDII）A"
$q3.Range("D19").Value = "0.25"
$q3.Range("E19").Value = "84.65"
$q3.Range("F19").Value = "3.13"
$q3.Range("G19").Value = "0.0078"
$q3.Range("H19").Value = 7

$q3.Range("B20").Value = "012871"
$q3.Range("C20").Value = "易方达纳斯达克100指数美元（QDII-LOF）C"
$q3.Range("D20").Value = "0.18"
$q3.Range("E20").Value = "90.67"
$q3.Range("F20").Value = "3.25"
$q3.Range("G20").Value = "0.0058"
$q3.Range("H20").Value = 5

$q3.Range("B21").Value = "012870"
$q3.Range("C21").Value = "易方达纳斯达克100指数人民币（QDII-LOF）C"
$q3.Range("D21").Value = "0.18"
$q3.Range("E21").Value = "90.67"
$q3.Range("F21").Value = "3.25"
$q3.Range("G21").Value = "0.0058"
$q3.Range("H21").Value = 5

$q3.Range("B22").Value = "012861"
$q3.Range("C22").Value = "易方达标普500指数（QDII-LOF）美元 C"
$q3.Range("D22").Value = "0.08"
$q3.Range("E22").Value = "90.72"
$q3.Range("F22").Value = "1.56"
$q3.Range("G22").Value = "0.0012"
$q3.Range("H22").Value = 6

$q3.Range("B23").Value = "014002"
$q3.Range("C23").Value = "浦银安盛全球智能科技股票（QDII）C"
$q3.Range("D23").Value = "0.01"
$q3.Range("E23").Value = "84.65"
$q3.Range("F23").Value = "3.13"
$q3.Range("G23").Value = "0.0003"
$q3.Range("H23").Value = 7

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: shift the quarterly rows down by one and
#    insert the brand-new 2022-Q3 row on top (copying formatting for the new
#    last row, which used to not exist).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 16
$total.Range("D7").Value = 7.71
$total.Range("A7").Value = 5

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 17
$total.Range("D6").Value = 9.890000000000001

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 19
$total.Range("D5").Value = 12.6

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 20
$total.Range("D4").Value = 21.42

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 22
$total.Range("D3").Value = 11.71

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 22
$total.Range("D2").Value = 10.29
